$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- new values (previously held by row 3)
$ws.Range("A2").Value = 111809656
$ws.Range("B2").Value = 56398
$ws.Range("E2").Value = 100109
$ws.Range("F2").Value = "Tretåig hackspett"
$ws.Range("G2").Value = "Picoides tridactylus"
$ws.Range("H2").Value = "(Linnaeus, 1758)"
$ws.Range("Q2").Value = 610542.5625081829
$ws.Range("R2").Value = 7180707.182562917
$ws.Range("Z2").Value = "14:37"
$ws.Range("AB2").Value = "14:37"

# Row 3 <- new values (previously held by row 2)
$ws.Range("A3").Value = 111809580
$ws.Range("B3").Value = 77515
$ws.Range("E3").Value = 6425
$ws.Range("F3").Value = "Garnlav"
$ws.Range("G3").Value = "Alectoria sarmentosa"
$ws.Range("H3").Value = "(Ach.) Ach."
$ws.Range("Q3").Value = 610571.4165256479
$ws.Range("R3").Value = 7180702.680798599
$ws.Range("Z3").Value = "14:31"
$ws.Range("AB3").Value = "14:31"

# Row 4 <- id and time updates
$ws.Range("A4").Value = 111809606
$ws.Range("Z4").Value = "14:33"
$ws.Range("AB4").Value = "14:33"
